$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-13 20:50:28"
$wsZhCn.Range("H4").Value = "2016-03-13 20:50:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-13 20:50:32"
$wsDeDe.Range("H4").Value = "2016-03-13 20:50:56"
